# Apply updated crypto price/volume data (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.633.56"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "2.169.75"
$ws.Range("E3").Value = "  -2.62%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'238.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("D7").Value = "'72.11"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.93%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -4.45%  "
$ws.Range("D10").Value = "'40.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.28%  "
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("D12").Value = "'54.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.100"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").Value = "2.496.91"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").Value = "'14.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "2.170.41"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("E18").Value = "  -6.85%  "
$ws.Range("D19").Value = "41.519.92"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "'70.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.00%  "
$ws.Range("D22").Value = "'5.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.01%  "
$ws.Range("D23").Value = "'9.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -14.23%  "
$ws.Range("D24").Value = "'226.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'10.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.97%  "
$ws.Range("E28").Value = "  -9.65%  "
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").Value = "'170.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'33.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.49%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'19.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("D34").Value = "'0.0771"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").Value = "'5.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.75%  "
$ws.Range("D36").Value = "'0.120"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.42%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -4.32%  "
$ws.Range("D39").Value = "'0.0306"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "'12.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.46%  "
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").Value = "'5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.57%  "
$ws.Range("D43").Value = "'59.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.17%  "
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("D46").Value = "'0.0965"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.65%  "
$ws.Range("D47").Value = "'97.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.54%  "
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("E49").Value = "  -4.62%  "
$ws.Range("E50").Value = "  -7.06%  "
$ws.Range("E51").Value = "  -2.34%  "

Write-Output "Updated 87 cells (24 numeric-looking forced to text, 63 plain text)."
